$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B7").Value = 0
$ws.Range("C7").Value = -0.05
$ws.Range("D7").Value = -0.02
$ws.Range("H7").Value = -0
$ws.Range("J7").Value = 0.02
$ws.Range("N7").Value = -0.05
